# Apply updated team-specific time-data matrix values (Akron_B) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2" = 0.2008196721311475;
    "C2" = 0.4672131147540984;
    "J2" = 0.01229508196721311;
    "P2" = 0.1926229508196721;
    "S2" = 0.1270491803278689;
    "C3" = 0.02542372881355932;
    "J3" = 0.03389830508474576;
    "P3" = 0.6949152542372882;
    "S3" = 0.2457627118644068;
    "J4" = 0.07692307692307693;
    "P4" = 0.6153846153846154;
    "S4" = 0.3076923076923077;
    "B6" = 0.05633802816901409;
    "D6" = 0.009389671361502348;
    "F6" = 0.05164319248826291;
    "J6" = 0.244131455399061;
    "O6" = 0.009389671361502348;
    "Q6" = 0.1126760563380282;
    "R6" = 0.07511737089201878;
    "S6" = 0.4413145539906103;
    "B7" = 0.09497206703910614;
    "D7" = 0.01675977653631285;
    "F7" = 0.0670391061452514;
    "J7" = 0.1340782122905028;
    "O7" = 0.0335195530726257;
    "Q7" = 0.1787709497206704;
    "R7" = 0.06145251396648044;
    "S7" = 0.4134078212290503;
    "B8" = 0.05489260143198091;
    "D8" = 0.01193317422434368;
    "F8" = 0.06443914081145585;
    "J8" = 0.1002386634844869;
    "O8" = 0.02147971360381861;
    "Q8" = 0.2124105011933174;
    "R8" = 0.09307875894988067;
    "S8" = 0.441527446300716;
    "B9" = 0.07906976744186046;
    "D9" = 0.01395348837209302;
    "F9" = 0.06046511627906977;
    "J9" = 0.1069767441860465;
    "O9" = 0.03255813953488372;
    "Q9" = 0.1953488372093023;
    "R9" = 0.06511627906976744;
    "S9" = 0.4465116279069767;
    "B10" = 0.1021159153633855;
    "D10" = 0.02391904323827047;
    "F10" = 0.08555657773689053;
    "J10" = 0.1131554737810488;
    "O10" = 0.01195952161913524;
    "Q10" = 0.2005519779208832;
    "R10" = 0.07727690892364306;
    "S10" = 0.3854645814167433;
    "G11" = 0.137546468401487;
    "J11" = 0.07434944237918216;
    "K11" = 0.1635687732342007;
    "L11" = 0.5985130111524164;
    "S11" = 0.02602230483271376;
    "G12" = 0.7559523809523809;
    "J12" = 0.1428571428571428;
    "K12" = 0.005952380952380952;
    "L12" = 0.03571428571428571;
    "S12" = 0.05952380952380952;
    "F13" = 0.025;
    "G13" = 0.575;
    "J13" = 0.35;
    "S13" = 0.05;
    "F15" = 0.009174311926605505;
    "H15" = 0.1467889908256881;
    "I15" = 0.07798165137614679;
    "J15" = 0.3669724770642202;
    "K15" = 0.05504587155963303;
    "M15" = 0.01376146788990826;
    "O15" = 0.05963302752293578;
    "S15" = 0.2706422018348624;
    "F16" = 0.05555555555555555;
    "H16" = 0.1805555555555556;
    "I16" = 0.09027777777777778;
    "J16" = 0.3333333333333333;
    "K16" = 0.125;
    "M16" = 0.02083333333333333;
    "O16" = 0.04166666666666666;
    "S16" = 0.1527777777777778;
    "F17" = 0.01;
    "H17" = 0.21;
    "I17" = 0.1;
    "J17" = 0.36;
    "K17" = 0.11;
    "M17" = 0.015;
    "N17" = 0.0025;
    "O17" = 0.055;
    "S17" = 0.1375;
    "F18" = 0.006211180124223602;
    "H18" = 0.1801242236024845;
    "I18" = 0.09316770186335403;
    "J18" = 0.4472049689440994;
    "K18" = 0.06832298136645963;
    "M18" = 0.01863354037267081;
    "O18" = 0.08695652173913043;
    "S18" = 0.09937888198757763;
    "F19" = 0.01529790660225443;
    "H19" = 0.2053140096618357;
    "I19" = 0.1046698872785829;
    "J19" = 0.3446054750402576;
    "K19" = 0.106280193236715;
    "M19" = 0.02012882447665056;
    "O19" = 0.08293075684380032;
    "S19" = 0.1207729468599034
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
